# Generate Report for Handoff
# Updates the handoff identifiers (old UUID -> new UUID) and the
# associated handoff/handback timestamps across the Overview, zh-cn and
# de-de worksheets, including the hyperlink display text that mirrors
# the cell values.

$wb = $excel.ActiveWorkbook

$oldId = "5f04f9b5-1c71-4b93-9000-9590c691d1fe"
$newId = "cd85353c-1e1a-46cf-a74c-14c80df945f9"
$oldHash = "482404404f7c5e61c5c3d200e535dbc8bf00dd98"
$newHash = "9b380ae25148ae2323a51d5721c430edf0bb2d6a"

$newMdName = "$newId.md"
$newZhXlfName = "$newId.$newHash.zh-cn.xlf"
$newDeXlfName = "$newId.$newHash.de-de.xlf"

$newOverviewDate = "2016-03-24 10:39:09"
$newXlfDate = "2016-03-24 10:38:59"

$mdAddr = "https://github.com/OpenLocalizationTest/oltest/blob/3b1fe54c1a372ed08397e4759519cc9cf90d6fad/e2e/$oldId.md"
$zhXlfAddr = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/35fb0c9ffc5a920ca7603a23b130e2d7bb938704/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/$oldId.$oldHash.zh-cn.xlf"
$deXlfAddr = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9c15b33e6f6acec90d95368f4d2ef1872955e66e/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/$oldId.$oldHash.de-de.xlf"

# ----------------------------------------------------------------------
# Sheet "Overview" - single hyperlinked cell (A2)
# ----------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# NOTE: Range.Hyperlinks.Delete() removes every hyperlink on the sheet,
# so all of a sheet's hyperlinks must be re-added together afterwards.
$wsOverview.Range("A2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $mdAddr, "", "", $newMdName) | Out-Null

$wsOverview.Range("D2").Value = $newOverviewDate

# ----------------------------------------------------------------------
# Sheet "zh-cn" - two hyperlinked cells (A2, D2)
# ----------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $mdAddr, "", "", $newMdName) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D2"), $zhXlfAddr, "", "", $newZhXlfName) | Out-Null

$wsZhCn.Range("E2").Value = $newXlfDate

# ----------------------------------------------------------------------
# Sheet "de-de" - two hyperlinked cells (A2, D2)
# ----------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $mdAddr, "", "", $newMdName) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D2"), $deXlfAddr, "", "", $newDeXlfName) | Out-Null

$wsDeDe.Range("E2").Value = $newOverviewDate

"Handoff report regenerated"
